# ajuste en contratos ws
#
# - Cambia el protocolo de "SOAP" a "REST" en ambas hojas (celda C6).
# - Agrega el campo "idUsuario" a la seccion "Datos entrada" de
#   ws_servidor_autorizarEntregarDulce (fila 19) y de
#   ws_servidor_recibirEntregaDulce (fila 13).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # autorizarEntregarDulce
$ws2 = $wb.Worksheets.Item(2)   # recibirEntregaDulce

# --- Protocolo: SOAP -> REST (ambas hojas comparten el mismo texto) ---
$ws1.Range("C6").Value = "REST"
$ws2.Range("C6").Value = "REST"

# --- Hoja "autorizarEntregarDulce": nueva fila de datos de entrada (fila 19) ---

# Copiamos el formato "en blanco" (idéntico al de F8 en la otra hoja) para
# las columnas Campo / Tipo Dato / Comentarios.
$ws2.Range("F8").Copy()
$ws1.Range("A19").PasteSpecial(-4122)
$ws1.Range("C19").PasteSpecial(-4122)
$ws1.Range("D19").PasteSpecial(-4122)

# Copiamos el formato normal (igual al resto de celdas de texto de la tabla)
# para la columna Descripcion.
$ws1.Range("A11").Copy()
$ws1.Range("B19").PasteSpecial(-4122)

$ws1.Range("A19").Value = "idUsuario"
$ws1.Range("B19").Value = "id del usuario que se le va a entregar el dulce"
$ws1.Range("C19").Value = "int"
$ws1.Range("D19").Value = "id del usuario que se le va a entregar el dulce"

# --- Hoja "recibirEntregaDulce": nueva fila de datos de entrada (fila 13) ---

$ws2.Range("F8").Copy()
$ws2.Range("A13").PasteSpecial(-4122)
$ws2.Range("C13").PasteSpecial(-4122)

$ws2.Range("A11").Copy()
$ws2.Range("B13").PasteSpecial(-4122)
$ws2.Range("D13").PasteSpecial(-4122)

$ws2.Range("A13").Value = "idUsuario"
$ws2.Range("B13").Value = "id del usuario que se le entrego el dulce"
$ws2.Range("C13").Value = "int"
$ws2.Range("D13").Value = "id del usuario que se le entrego el dulce"

$excel.CutCopyMode = $false
